# Update the "cryptos" worksheet with refreshed price/volume data.
# Numeric-looking text in column D is prefixed with a leading apostrophe so
# Excel keeps it as text (matching the source data's inlineStr type) instead
# of auto-converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "56.965.97"
$ws.Range('E2').Value = "  -2.02%  "
$ws.Range('D3').Value = "3.065.94"
$ws.Range('E3').Value = "  -1.61%  "
$ws.Range('E4').Value = "  +0.03%  "
$ws.Range('D5').Value = "'520.54"
$ws.Range('E5').Value = "  -1.28%  "
$ws.Range('D6').Value = "'134.89"
$ws.Range('E6').Value = "  -5.16%  "
$ws.Range('E7').Value = "  +0.02%  "
$ws.Range('D8').Value = "3.065.88"
$ws.Range('E8').Value = "  -1.57%  "
$ws.Range('E9').Value = "  +2.19%  "
$ws.Range('D10').Value = "'7.32"
$ws.Range('E11').Value = "  -2.85%  "
$ws.Range('E12').Value = "  +0.79%  "
$ws.Range('B13').Value = "TRON"
$ws.Range('C13').Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range('D13').Value = "'0.135"
$ws.Range('E13').Value = "  +0.98%  "
$ws.Range('B14').Value = "WrappedliquidstakedEther2.0"
$ws.Range('C14').Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D14').Value = "3.600.35"
$ws.Range('E14').Value = "  -1.40%  "
$ws.Range('D15').Value = "'25.17"
$ws.Range('E15').Value = "  -2.51%  "
$ws.Range('E16').Value = "  -3.04%  "
$ws.Range('D17').Value = "57.004.25"
$ws.Range('E17').Value = "  -2.03%  "
$ws.Range('D18').Value = "3.077.32"
$ws.Range('E18').Value = "  -1.11%  "
$ws.Range('E19').Value = "  -4.53%  "
$ws.Range('D20').Value = "'12.40"
$ws.Range('E20').Value = "  -2.98%  "
$ws.Range('D21').Value = "'7.80"
$ws.Range('E21').Value = "  -2.51%  "
$ws.Range('D22').Value = "'347.32"
$ws.Range('E22').Value = "  +1.14%  "
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = "  -0.10%  "
$ws.Range('D24').Value = "'68.52"
$ws.Range('E24').Value = "  +1.28%  "
$ws.Range('D25').Value = "'0.496"
$ws.Range('E25').Value = "  -3.87%  "
$ws.Range('E26').Value = "  -2.23%  "
$ws.Range('E27').Value = "  +0.05%  "
$ws.Range('D28').Value = "0.0₃0863"
$ws.Range('E28').Value = "  -7.13%  "
$ws.Range('E29').Value = "  -0.20%  "
$ws.Range('D30').Value = "'7.17"
$ws.Range('E30').Value = "  -1.54%  "
$ws.Range('E31').Value = "  -1.27%  "
$ws.Range('D32').Value = "'5.83"
$ws.Range('E32').Value = "  -8.63%  "
$ws.Range('D33').Value = "'20.84"
$ws.Range('B34').Value = "Monero"
$ws.Range('C34').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D34').Value = "'159.24"
$ws.Range('E34').Value = "  +0.38%  "
$ws.Range('B35').Value = "NEARProtocol"
$ws.Range('C35').Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('D35').Value = "'4.82"
$ws.Range('E35').Value = "  +3.59%  "
$ws.Range('E36').Value = "  -5.75%  "
$ws.Range('D37').Value = "'5.98"
$ws.Range('E37').Value = "  -3.56%  "
$ws.Range('D38').Value = "'25.38"
$ws.Range('E38').Value = "  -3.99%  "
$ws.Range('E39').Value = "  -1.87%  "
$ws.Range('D40').Value = "'0.0656"
$ws.Range('E40').Value = "  -2.03%  "
$ws.Range('E41').Value = "  -4.14%  "
$ws.Range('D42').Value = "'4.01"
$ws.Range('E42').Value = "  +0.27%  "
$ws.Range('D43').Value = "'0.690"
$ws.Range('E43').Value = "  +0.09%  "
$ws.Range('D44').Value = "2.388.00"
$ws.Range('E44').Value = "  +5.00%  "
$ws.Range('D45').Value = "'36.61"
$ws.Range('E45').Value = "  +0.31%  "
$ws.Range('E46').Value = "  +0.06%  "
$ws.Range('D47').Value = "3.108.25"
$ws.Range('E47').Value = "  -1.55%  "
$ws.Range('D48').Value = "'0.0261"
$ws.Range('E48').Value = "  -0.77%  "
$ws.Range('D49').Value = "'0.950"
$ws.Range('E49').Value = "  -4.90%  "
$ws.Range('E50').Value = "  -2.98%  "
$ws.Range('D51').Value = "'19.59"
$ws.Range('E51').Value = "  -5.07%  "
